# Opponent_List.xlsx — position corrections on the "Positions" lookup sheet,
# plus turning that range into a proper AutoFilter table.

$wb = $excel.ActiveWorkbook

$posSheet = $wb.Worksheets.Item("Positions")
$oppSheet = $wb.Worksheets.Item("Opponent_List")

# --- Correct mis-entered positions (column B) -----------------------------
$posSheet.Range("B3").Value = "SG"
$posSheet.Range("B8").Value = "PG"
$posSheet.Range("B16").Value = "SF"
$posSheet.Range("B24").Value = "SF"
$posSheet.Range("B36").Value = "C"
$posSheet.Range("B52").Value = "C"
$posSheet.Range("B65").Value = "SG"
$posSheet.Range("B69").Value = "C"
$posSheet.Range("B93").Value = "C"
$posSheet.Range("B102").Value = "PF"
$posSheet.Range("B103").Value = "SG"
$posSheet.Range("B106").Value = "PF"
$posSheet.Range("B122").Value = "SF"
$posSheet.Range("B145").Value = "SG"
$posSheet.Range("B162").Value = "PG"
$posSheet.Range("B193").Value = "C"
$posSheet.Range("B230").Value = "SF"
$posSheet.Range("B251").Value = "PF"
$posSheet.Range("B259").Value = "PG"
$posSheet.Range("B267").Value = "PF"
$posSheet.Range("B281").Value = "SF"
$posSheet.Range("B288").Value = "C"
$posSheet.Range("B310").Value = "SF"
$posSheet.Range("B315").Value = "SF"
$posSheet.Range("B348").Value = "C"
$posSheet.Range("B355").Value = "C"
$posSheet.Range("B386").Value = "SF"
$posSheet.Range("B391").Value = "SG"
$posSheet.Range("B395").Value = "SG"
$posSheet.Range("B417").Value = "SG"
$posSheet.Range("B419").Value = "PF"
$posSheet.Range("B451").Value = "SF"
$posSheet.Range("B479").Value = "SG"
$posSheet.Range("B505").Value = "SF"
$posSheet.Range("B547").Value = "PG"
$posSheet.Range("B550").Value = "SF"
$posSheet.Range("B562").Value = "SG"
$posSheet.Range("B572").Value = "PF"
$posSheet.Range("B576").Value = "PG"
$posSheet.Range("B589").Value = "PF"
$posSheet.Range("B596").Value = "C"

# --- Turn the lookup range into an AutoFilter table -----------------------
$null = $posSheet.Range("A1:B603").AutoFilter()
$filterName = $posSheet.Names.Add("_xlnm._FilterDatabase", "=Positions!`$A`$1:`$B`$603")
$filterName.Visible = $false

# --- Restore cursor positions on each sheet --------------------------------
$null = $posSheet.Range("B6").Select()

$null = $oppSheet.Activate()
$null = $oppSheet.Range("R17").Select()
